# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 11:52"

# --- Australia (row 47) ---
$ws.Range("B47").Value = 6753
$ws.Range("C47").Value = 7
$ws.Range("E47").Value = 947
$ws.Range("F47").Value = 34

# --- Moldavia (row 58) ---
$ws.Range("D58").Value = 1182
$ws.Range("E58").Value = 2473
$ws.Range("F58").Value = 237
$ws.Range("G58").Value = 5
$ws.Range("H58").Value = 116

# --- Reordering of Armenia..Bulgaria block (rows 68-78) ---
# A new country (Ghana) is inserted before Armenia, and Estonia is
# inserted before Bosnia y Herzegovina; the old Ghana/Estonia rows are
# removed from their former spot. Net effect on the fixed row positions
# 68-78 is a relabel + new stats as below.

# Row 68: now Ghana (new data)
$ws.Range("A68").Value = "Ghana"
$ws.Range("B68").Value = 2074
$ws.Range("C68").Value = 403
$ws.Range("D68").Value = 212
$ws.Range("E68").Value = 1845
$ws.Range("F68").Value = 4
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 17

# Row 69: now Armenia
$ws.Range("A69").Value = "Armenia"
$ws.Range("B69").Value = 2066
$ws.Range("C69").Value = 134
$ws.Range("D69").Value = 929
$ws.Range("E69").Value = 1105
$ws.Range("F69").Value = 10
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 32

# Row 70: now Croacia
$ws.Range("A70").Value = "Croacia"
$ws.Range("B70").Value = 2062
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 1288
$ws.Range("E70").Value = 707
$ws.Range("F70").Value = 19
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 67

# Row 71: now Uzbekistan
$ws.Range("A71").Value = "Uzbekistan"
$ws.Range("B71").Value = 2017
$ws.Range("C71").Value = 15
$ws.Range("D71").Value = 1096
$ws.Range("E71").Value = 912
$ws.Range("F71").Value = 8
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 9

# Row 72: now Irak
$ws.Range("A72").Value = "Irak"
$ws.Range("B72").Value = 2003
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 1346
$ws.Range("E72").Value = 565
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 92

# Row 73: now Camerun
$ws.Range("A73").Value = "Camerun"
$ws.Range("B73").Value = 1832
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 934
$ws.Range("E73").Value = 837
$ws.Range("F73").Value = 12
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 61

# Row 74: now Islandia
$ws.Range("A74").Value = "Islandia"
$ws.Range("B74").Value = 1797
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 1656
$ws.Range("E74").Value = 131
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 10

# Row 75: now Azerbaiyan
$ws.Range("A75").Value = "Azerbaiyan"
$ws.Range("B75").Value = 1766
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 1267
$ws.Range("E75").Value = 476
$ws.Range("F75").Value = 18
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 23

# Row 76: now Nigeria
$ws.Range("A76").Value = "Nigeria"
$ws.Range("B76").Value = 1728
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 307
$ws.Range("E76").Value = 1370
$ws.Range("F76").Value = 2
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 51

# Row 77: now Estonia (new data)
$ws.Range("A77").Value = "Estonia"
$ws.Range("B77").Value = 1689
$ws.Range("C77").Value = 23
$ws.Range("D77").Value = 249
$ws.Range("E77").Value = 1388
$ws.Range("F77").Value = 9
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 52

# Row 78: now Bosnia y Herzegovina
$ws.Range("A78").Value = "Bosnia y Herzegovina"
$ws.Range("B78").Value = 1677
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 710
$ws.Range("E78").Value = 902
$ws.Range("F78").Value = 4
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 65

# --- Etiopia (row 142) ---
$ws.Range("B142").Value = 131
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 59

# --- Sudan del Sur (row 178) ---
$ws.Range("B178").Value = 35
$ws.Range("C178").Value = 1
$ws.Range("E178").Value = 35
